$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Copy the format of an existing fully-populated date row (113, style s=2)
# onto the two brand-new trailing rows (114, 115) before writing their values,
# so the new cells inherit the same cellXf (bold border date style for col A,
# plain numeric for B/C/D) instead of the blank default style.
$ws.Range("A113:D113").Copy()
$ws.Range("A114:D114").PasteSpecial(-4122)
$ws.Range("A113:D113").Copy()
$ws.Range("A115:D115").PasteSpecial(-4122)

# Write the updated data block (rows 90-115): a row for 2021-02-08 (serial 44235)
# was inserted, shifting the previously-unpublished trailing rows down by one and
# adding one further new trailing row, plus refreshed rolling-sum figures.
$ws.Range("A90").Value = 44232
$ws.Range("B90").Value = 2
$ws.Range("C90").Value = 17
$ws.Range("D90").Value = 112.679790548154

$ws.Range("A91").Value = 44233
$ws.Range("B91").Value = 2
$ws.Range("C91").Value = 13
$ws.Range("D91").Value = 86.16689865447074

$ws.Range("A92").Value = 44234
$ws.Range("B92").Value = 6
$ws.Range("C92").Value = 12
$ws.Range("D92").Value = 79.53867568104991

$ws.Range("A93").Value = 44235
$ws.Range("B93").Value = 2
$ws.Range("C93").Value = 15
$ws.Range("D93").Value = 99.42334460131239

$ws.Range("A94").Value = 44236
$ws.Range("B94").Value = 0
$ws.Range("C94").Value = 18
$ws.Range("D94").Value = 119.3080135215749

$ws.Range("A95").Value = 44237
$ws.Range("B95").Value = 0
$ws.Range("C95").Value = 19
$ws.Range("D95").Value = 125.9362364949957

$ws.Range("A96").Value = 44238
$ws.Range("B96").Value = 3
$ws.Range("C96").Value = 14
$ws.Range("D96").Value = 92.79512162789156

$ws.Range("A97").Value = 44239
$ws.Range("B97").Value = 5
$ws.Range("C97").Value = 18
$ws.Range("D97").Value = 119.3080135215749

$ws.Range("A98").Value = 44240
$ws.Range("B98").Value = 3
$ws.Range("C98").Value = 21
$ws.Range("D98").Value = 139.1926824418373

$ws.Range("A99").Value = 44241
$ws.Range("B99").Value = 1
$ws.Range("C99").Value = 22
$ws.Range("D99").Value = 145.8209054152582

$ws.Range("A100").Value = 44242
$ws.Range("B100").Value = 6
$ws.Range("C100").Value = 27
$ws.Range("D100").Value = 178.9620202823623

$ws.Range("A101").Value = 44243
$ws.Range("B101").Value = 3
$ws.Range("C101").Value = 29
$ws.Range("D101").Value = 192.218466229204

$ws.Range("A102").Value = 44244
$ws.Range("B102").Value = 1
$ws.Range("C102").Value = 36
$ws.Range("D102").Value = 238.6160270431498

$ws.Range("A103").Value = 44245
$ws.Range("B103").Value = 8
$ws.Range("C103").Value = 39
$ws.Range("D103").Value = 258.5006959634122

$ws.Range("A104").Value = 44246
$ws.Range("B104").Value = 7
$ws.Range("C104").Value = 43
$ws.Range("D104").Value = 285.0135878570956

$ws.Range("A105").Value = 44247
$ws.Range("B105").Value = 10
$ws.Range("C105").Value = 48
$ws.Range("D105").Value = 318.1547027241996

$ws.Range("A106").Value = 44248
$ws.Range("B106").Value = 4
$ws.Range("C106").Value = 50
$ws.Range("D106").Value = 331.4111486710413

$ws.Range("A107").Value = 44249
$ws.Range("B107").Value = 10
$ws.Range("C107").Value = 56
$ws.Range("D107").Value = 371.1804865115662

$ws.Range("A108").Value = 44250
$ws.Range("B108").Value = 8
$ws.Range("C108").Value = 66
$ws.Range("D108").Value = 437.4627162457745

$ws.Range("A109").Value = 44251
$ws.Range("B109").Value = 3
$ws.Range("C109").Value = 72
$ws.Range("D109").Value = 477.2320540862995

$ws.Range("A110").Value = 44252
$ws.Range("B110").Value = 14
$ws.Range("C110").Value = 78
$ws.Range("D110").Value = 517.0013919268245

$ws.Range("A111").Value = 44253
$ws.Range("B111").Value = 17
$ws.Range("C111").Value = 102
$ws.Range("D111").Value = 676.0787432889242

$ws.Range("A112").Value = 44254
$ws.Range("B112").Value = 16
$ws.Range("C112").Value = 100
$ws.Range("D112").Value = 662.8222973420826

$ws.Range("A113").Value = 44255
$ws.Range("B113").Value = 10
$ws.Range("C113").Value = ""
$ws.Range("D113").Value = ""

$ws.Range("A114").Value = 44256
$ws.Range("B114").Value = 34
$ws.Range("C114").Value = ""
$ws.Range("D114").Value = ""

$ws.Range("A115").Value = 44257
$ws.Range("B115").Value = 6
$ws.Range("C115").Value = ""
$ws.Range("D115").Value = ""
